$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "cheaper / smaller rect" part tracking row (row 8):
#   A8 = Part, B8 = Description, C8 = specific product link, D8 = filter-search link
$ws.Range("B8").Value = "cheaper / smaller rect"
$ws.Range("A8").Value = "BAV199S"
$ws.Range("D8").Value = "https://www.digikey.de/en/products/filter/diodes-rectifiers-arrays/286?s=N4IgjCBcoGwJxVAYygMwIYBsDOBTANCAPZQDaIALAAxwDMdIh1cFtVIAuoQA4AuUIAMq8ATgEsAdgHMQAX0JwA7BGggUkDDgLEy4AExg4ADj2ce-SENGSZ8kAFpTq9aICu2kpHIBWRiAQcsnYwiCBiACYC9mBUEOYCfrwAnty4AujYKEFAA"
$ws.Range("C8").Value = "https://www.digikey.de/en/products/detail/panjit-international-inc/BAV199S-R1-00001/14661140"

# Move / reflect the active selection on the sheet to C8
$ws.Range("C8").Select()
